# Delete the data row for "「子供に食べさせている忍耐強い母親」" (previously row 668).
# This removes the whole row and shifts all following rows up by one,
# matching the diff where rows 669-724 become 668-723 and the sheet
# dimension shrinks from A1:C724 to A1:C723.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(668).Delete()
